# Update "想去人数" (want-to-go count, column F) figures across the three
# sheets that carry per-event rows: 展览 (sheet 1), 演出 (sheet 2) and
# 全部类型 (sheet 4). 本地生活 (sheet 3) only has a header row, so it is
# untouched. Each row keeps the same event (same row index) but its
# interest-count value is refreshed to the newer scrape.

$wb = $excel.ActiveWorkbook

# --- 展览 (Worksheet 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value  = 817
$ws.Range("F5").Value  = 1157
$ws.Range("F6").Value  = 3370
$ws.Range("F7").Value  = 2482
$ws.Range("F9").Value  = 2304
$ws.Range("F11").Value = 1945
$ws.Range("F14").Value = 619
$ws.Range("F15").Value = 83
$ws.Range("F16").Value = 287
$ws.Range("F17").Value = 17
$ws.Range("F18").Value = 40
$ws.Range("F19").Value = 257
$ws.Range("F20").Value = 60
$ws.Range("F21").Value = 411
$ws.Range("F23").Value = 84
$ws.Range("F24").Value = 450
$ws.Range("F25").Value = 647
$ws.Range("F26").Value = 70
$ws.Range("F28").Value = 341
$ws.Range("F29").Value = 27
$ws.Range("F30").Value = 1596
$ws.Range("F31").Value = 732
$ws.Range("F32").Value = 738
$ws.Range("F33").Value = 1873
$ws.Range("F34").Value = 200
$ws.Range("F35").Value = 488
$ws.Range("F37").Value = 528
$ws.Range("F38").Value = 1161
$ws.Range("F40").Value = 392

# --- 演出 (Worksheet 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 55
$ws.Range("F8").Value = 4

# --- 全部类型 (Worksheet 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value  = 817
$ws.Range("F5").Value  = 1157
$ws.Range("F6").Value  = 3370
$ws.Range("F7").Value  = 2482
$ws.Range("F9").Value  = 2304
$ws.Range("F11").Value = 1945
$ws.Range("F14").Value = 619
$ws.Range("F15").Value = 83
$ws.Range("F16").Value = 287
$ws.Range("F17").Value = 17
$ws.Range("F18").Value = 40
$ws.Range("F19").Value = 257
$ws.Range("F20").Value = 60
$ws.Range("F21").Value = 411
$ws.Range("F23").Value = 84
$ws.Range("F24").Value = 450
$ws.Range("F25").Value = 647
$ws.Range("F26").Value = 70
$ws.Range("F27").Value = 55
$ws.Range("F31").Value = 341
$ws.Range("F32").Value = 27
$ws.Range("F33").Value = 1596
$ws.Range("F34").Value = 732
$ws.Range("F36").Value = 738
$ws.Range("F37").Value = 1873
$ws.Range("F38").Value = 200
$ws.Range("F41").Value = 4
$ws.Range("F42").Value = 488
$ws.Range("F44").Value = 528
$ws.Range("F45").Value = 1161
$ws.Range("F47").Value = 392
